$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new condition row (row 35): cityHasTasks / 城市是否有任务 / city / taskNumber / > / number / 0
$ws.Range("A35").Value = "cityHasTasks"
$ws.Range("B35").Value = "城市是否有任务"
$ws.Range("C35").Value = "city"
$ws.Range("D35").Value = "taskNumber"
$ws.Range("E35").Value = ">"
$ws.Range("F35").Value = "number"
$ws.Range("G35").Value = 0

# Match formatting used by the rest of column B (Chinese description text uses
# the 宋体 font style already defined in the workbook) by copying the format
# from an existing description cell instead of creating a brand-new style.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null

# Move the selection to the newly added last cell, as in the saved workbook.
$ws.Range("G35").Select() | Out-Null
